# Remove obsolete finite-difference-algorithm translation rows from the
# fr-FR localization table. The corresponding English source strings
# (strRadBackwardOne, strRadCentralFive, strRadCentralThree,
# strRadForwardOne) were removed upstream, so their rows are dropped here
# too. The sheet's Table1/"Tabla13" data is sorted alphabetically by Key
# (column B), so deleting these rows shifts everything below them up by
# one row each -- which is exactly the shape of the target edit (table
# and sheet dimension shrink from B2:E169 to B2:E165).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fr-FR")

# Keys (column B) of the rows to remove, in the order they currently
# appear top-to-bottom on the sheet.
$keysToRemove = @(
    "strRadBackwardOne",
    "strRadCentralFive",
    "strRadCentralThree",
    "strRadForwardOne"
)

# Find the current row number for each key, then delete bottom-up so
# already-found row numbers for the remaining keys stay valid while we
# work.
$rowsToDelete = @()
foreach ($key in $keysToRemove) {
    $found = $ws.Columns.Item(2).Find($key)
    $rowsToDelete += $found.Row
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
